$wb = $excel.ActiveWorkbook
$actor = $wb.Worksheets.Item("Actor")
$item = $wb.Worksheets.Item("Item")

# Actor sheet: fill in the empty "D" (desc) type cells
$actor.Range("D4").Value = "string"
$actor.Range("D5").Value = "desc"

# Item sheet: row 4 currently holds the field-name row and row 5 the
# type row (opposite order from the Actor sheet). Swap them so row 4 is
# the type row and row 5 is the field-name row, matching Actor.

# 1) swap the cell VALUES between row 4 and row 5
$row4 = @()
$row5 = @()
for ($c = 1; $c -le 9; $c++) {
    $row4 += $item.Cells.Item(4, $c).Value()
    $row5 += $item.Cells.Item(5, $c).Value()
}
for ($c = 1; $c -le 9; $c++) {
    $item.Cells.Item(4, $c).Value = $row5[$c - 1]
    $item.Cells.Item(5, $c).Value = $row4[$c - 1]
}

# 2) swap the cell FORMATTING (fill/border styles) between row 4 and row
#    5 too, using a scratch row (20) as a holding area and bounded
#    ranges (A:I) so the copy doesn't spill formatting across the whole row.
$item.Range("A4:I4").Copy()
$item.Range("A20:I20").PasteSpecial(-4122)

$item.Range("A5:I5").Copy()
$item.Range("A4:I4").PasteSpecial(-4122)

$item.Range("A20:I20").Copy()
$item.Range("A5:I5").PasteSpecial(-4122)

$item.Range("A20:I20").Clear()

# 3) fill in the previously-empty "D" (desc) type cells for both rows
$item.Range("D4").Value = "string"
$item.Range("D5").Value = "desc"

# Leave the selection/active sheet where the edits were made.
[void]$actor.Activate()
[void]$actor.Range("D6").Select()

[void]$item.Activate()
[void]$item.Range("D5").Select()

Write-Host "Done"
